# Auto-generated edit script: updates crypto price/volume table cells
# per the commit 'Updated cryptos list on Sun Oct 13 16:58:35 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to Text format first so they store as plain strings (matches the
# original inline-string cell type for column D).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values (row order matches the sheet).
$ws.Range("D2").Value = '62.673.32'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '2.453.97'
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '571.13'
$ws.Range("E5").Value = '  -0.95%  '
$ws.Range("D6").Value = '146.01'
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '0.530'
$ws.Range("E8").Value = '  -1.77%  '
$ws.Range("D9").Value = '0.111'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("E11").Value = '  -1.44%  '
$ws.Range("E12").Value = '  -1.80%  '
$ws.Range("D13").Value = '28.59'
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("E14").Value = '  -3.01%  '
$ws.Range("D15").Value = '2.900.22'
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '62.471.89'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").Value = '2.451.58'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = '7.60'
$ws.Range("E18").Value = '  -6.13%  '
$ws.Range("D19").Value = '10.73'
$ws.Range("E19").Value = '  -3.14%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '321.30'
$ws.Range("E20").Value = '  -2.78%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '4.13'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = '2.20'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '9.89'
$ws.Range("E24").Value = '  +4.44%  '
$ws.Range("D25").Value = '65.00'
$ws.Range("E25").Value = '  -2.14%  '
$ws.Range("D26").Value = '641.74'
$ws.Range("E26").Value = '  -3.76%  '
$ws.Range("D28").Value = '0.0₃0957'
$ws.Range("E28").Value = '  -3.72%  '
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -2.96%  '
$ws.Range("D31").Value = '7.81'
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("E32").Value = '  -3.23%  '
$ws.Range("E33").Value = '  -2.44%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  -3.71%  '
$ws.Range("D36").Value = '4.63'
$ws.Range("E36").Value = '  -3.26%  '
$ws.Range("D37").Value = '150.96'
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").Value = '0.365'
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").Value = '18.47'
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = '5.32'
$ws.Range("E40").Value = '  -3.47%  '
$ws.Range("D41").Value = '2.72'
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("D43").Value = '0.0₆0309'
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("D45").Value = '152.75'
$ws.Range("E45").Value = '  +2.69%  '
$ws.Range("D46").Value = '15.40'
$ws.Range("E46").Value = '  +1.55%  '
$ws.Range("E47").Value = '  -2.44%  '
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").Value = '20.14'
$ws.Range("E49").Value = '  -3.51%  '
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("E51").Value = '  -1.42%  '
